$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-20
# from serial 45174 (2023-09-05) to serial 45175 (2023-09-06)
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
